$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 20 (ALC)
$ws.Range("H20").Value = 2849.8333
$ws.Range("I20").Value = 1019.8
$ws.Range("K20").Value = 1019.8
$ws.Range("M20").Value = -789.8
# row 35 (ALC)
$ws.Range("H35").Value = 2849.8333
$ws.Range("I35").Value = 1019.8
$ws.Range("K35").Value = 1019.8
$ws.Range("M35").Value = -640.8
# row 64 (ALC)
$ws.Range("H64").Value = 252475
$ws.Range("I64").Value = 501750
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 501750
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -501502
$ws.Range("N64").Value = -3696
# row 67 (ALC)
$ws.Range("H67").Value = 252475
$ws.Range("I67").Value = 501750
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 501750
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -500892
$ws.Range("N67").Value = -4916
# row 116 (ALC)
$ws.Range("H116").Value = 2381.25
$ws.Range("I116").Value = 2222
$ws.Range("K116").Value = 2222
$ws.Range("M116").Value = 1220
# row 129 (ALC)
$ws.Range("H129").Value = 962.7727
$ws.Range("J129").Value = 1078.8235
$ws.Range("L129").Value = 3236.4705
$ws.Range("N129").Value = -13236.4705
# row 137 (ALC)
$ws.Range("H137").Value = 1925.0526
$ws.Range("I137").Value = 1476.4445
$ws.Range("K137").Value = 4429.333500000001
$ws.Range("M137").Value = -1879.333500000001
# row 138 (ALC)
$ws.Range("H138").Value = 8254.536
$ws.Range("I138").Value = 1808
$ws.Range("J138").Value = 14394.096
$ws.Range("K138").Value = 5424
$ws.Range("L138").Value = 43182.288
$ws.Range("M138").Value = -284
$ws.Range("N138").Value = -53462.288

$ws = $wb.Worksheets.Item("ARM")
# row 2 (ARM)
$ws.Range("H2").Value = 144213.14
$ws.Range("I2").Value = 1398.4
$ws.Range("J2").Value = 501250
$ws.Range("K2").Value = 1398.4
$ws.Range("L2").Value = 501250
$ws.Range("M2").Value = -1285.4
$ws.Range("N2").Value = -501476
# row 6 (ARM)
$ws.Range("H6").Value = 23268
$ws.Range("I6").Value = 25002
$ws.Range("J6").Value = 19800
$ws.Range("K6").Value = 25002
$ws.Range("L6").Value = 19800
$ws.Range("M6").Value = -24829
$ws.Range("N6").Value = -20146
# row 8 (ARM)
$ws.Range("H8").Value = 49800
$ws.Range("J8").Value = 49800
$ws.Range("L8").Value = 49800
$ws.Range("N8").Value = -50088
# row 32 (ARM)
$ws.Range("H32").Value = 28737.377
$ws.Range("I32").Value = 4745.091
$ws.Range("J32").Value = 248666.67
$ws.Range("K32").Value = 4745.091
$ws.Range("L32").Value = 248666.67
$ws.Range("M32").Value = -4458.091
$ws.Range("N32").Value = -249240.67
# row 55 (ARM)
$ws.Range("H55").Value = 12309.9
$ws.Range("J55").Value = 12455.444
$ws.Range("L55").Value = 12455.444
$ws.Range("N55").Value = -13085.444
# row 61 (ARM)
$ws.Range("H61").Value = 2249.457
$ws.Range("I61").Value = 1408.8125
$ws.Range("J61").Value = 2957.3684
$ws.Range("K61").Value = 1408.8125
$ws.Range("L61").Value = 2957.3684
$ws.Range("M61").Value = -1196.8125
$ws.Range("N61").Value = -3381.3684
# row 88 (ARM)
$ws.Range("H88").Value = 3600
$ws.Range("I88").Value = 3800
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 3800
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = -3394
$ws.Range("N88").Value = -4312
# row 91 (ARM)
$ws.Range("H91").Value = 3600
$ws.Range("I91").Value = 3800
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 3800
$ws.Range("L91").Value = 3500
$ws.Range("M91").Value = -2396
$ws.Range("N91").Value = -6308
# row 116 (ARM)
$ws.Range("H116").Value = 144213.14
$ws.Range("I116").Value = 1398.4
$ws.Range("J116").Value = 501250
$ws.Range("K116").Value = 1398.4
$ws.Range("L116").Value = 501250
$ws.Range("M116").Value = 895.5999999999999
$ws.Range("N116").Value = -505838
# row 136 (ARM)
$ws.Range("H136").Value = 2249.457
$ws.Range("I136").Value = 1408.8125
$ws.Range("J136").Value = 2957.3684
$ws.Range("K136").Value = 4226.4375
$ws.Range("L136").Value = 8872.1052
$ws.Range("M136").Value = -1676.4375
$ws.Range("N136").Value = -13972.1052

$ws = $wb.Worksheets.Item("BSM")
# row 3 (BSM)
$ws.Range("H3").Value = 144213.14
$ws.Range("I3").Value = 1398.4
$ws.Range("J3").Value = 501250
$ws.Range("K3").Value = 1398.4
$ws.Range("L3").Value = 501250
$ws.Range("M3").Value = -1284.4
$ws.Range("N3").Value = -501478
# row 16 (BSM)
$ws.Range("H16").Value = 9009
$ws.Range("J16").Value = 9009
$ws.Range("L16").Value = 9009
$ws.Range("N16").Value = -9349

$ws = $wb.Worksheets.Item("CRP")
# row 31 (CRP)
$ws.Range("H31").Value = 39328.168
$ws.Range("I31").Value = 1320.7142
$ws.Range("J31").Value = 58331.895
$ws.Range("K31").Value = 1320.7142
$ws.Range("L31").Value = 58331.895
$ws.Range("M31").Value = -1025.7142
$ws.Range("N31").Value = -58921.895
# row 34 (CRP)
$ws.Range("H34").Value = 39328.168
$ws.Range("I34").Value = 1320.7142
$ws.Range("J34").Value = 58331.895
$ws.Range("K34").Value = 1320.7142
$ws.Range("L34").Value = 58331.895
$ws.Range("M34").Value = -1118.7142
$ws.Range("N34").Value = -58735.895
# row 37 (CRP)
$ws.Range("H37").Value = 39900
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 39900
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 39900
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -40114
# row 39 (CRP)
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -1609
$ws.Range("N39").ClearContents()
# row 49 (CRP)
$ws.Range("H49").Value = 2000
$ws.Range("I49").Value = 2000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -1818
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# row 132 (GSM)
$ws.Range("H132").Value = 2612.9688
$ws.Range("I132").Value = 1592.4783
$ws.Range("J132").Value = 5220.8887
$ws.Range("K132").Value = 4777.4349
$ws.Range("L132").Value = 15662.6661
$ws.Range("M132").Value = -2247.4349
$ws.Range("N132").Value = -20722.6661

$ws = $wb.Worksheets.Item("LTW")
# row 24 (LTW)
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
# row 68 (LTW)
$ws.Range("H68").Value = 2634.0527
$ws.Range("I68").Value = 1366.6666
$ws.Range("J68").Value = 3774.7
$ws.Range("K68").Value = 1366.6666
$ws.Range("L68").Value = 3774.7
$ws.Range("M68").Value = -617.6666
$ws.Range("N68").Value = -5272.7
# row 71 (LTW)
$ws.Range("H71").Value = 2634.0527
$ws.Range("I71").Value = 1366.6666
$ws.Range("J71").Value = 3774.7
$ws.Range("K71").Value = 6833.333000000001
$ws.Range("L71").Value = 18873.5
$ws.Range("M71").Value = -3089.333000000001
$ws.Range("N71").Value = -26361.5
# row 132 (LTW)
$ws.Range("H132").Value = 2854.9268
$ws.Range("I132").Value = 2881.4
$ws.Range("J132").Value = 2700.5
$ws.Range("K132").Value = 8644.200000000001
$ws.Range("L132").Value = 8101.5
$ws.Range("M132").Value = -6114.200000000001
$ws.Range("N132").Value = -13161.5
# row 140 (LTW)
$ws.Range("H140").Value = 73409.664
$ws.Range("J140").Value = 73409.664
$ws.Range("L140").Value = 73409.664
$ws.Range("N140").Value = -83769.664

$ws = $wb.Worksheets.Item("WVR")
# row 24 (WVR)
$ws.Range("H24").Value = 514900
$ws.Range("J24").Value = 514900
$ws.Range("L24").Value = 514900
$ws.Range("N24").Value = -515360
# row 132 (WVR)
$ws.Range("H132").Value = 2239.157
$ws.Range("I132").Value = 2120.7805
$ws.Range("J132").Value = 2724.5
$ws.Range("K132").Value = 6362.3415
$ws.Range("L132").Value = 8173.5
$ws.Range("M132").Value = -3832.3415
$ws.Range("N132").Value = -13233.5
